$d = $word.ActiveDocument

$d.Content.Find.Execute("You’re invited to participate in an evaluation study of CrisisText: A chatbot developed by Parenting for Lifelong Health (PLH), World Vision (WV), and IDEMS International to strengthen parents, caregivers, and children. We’re doing this study to learn about your experience with CrisisText to make sure the chatbot helps families like yours.", $true, $false, $false, $false, $false, $true, 1, $false, "شما برای اشتراک در یک مطالعهٔ برای ارزیابی در مورد پیام رسان بحران «CrisisText» دعوت شده‌اید: یک چت‌بات که توسط برنامهٔ پرورش والدین برای زندگی سالم(PLH)،مؤسسهٔ ورلدویژن (WV)، و IDEMS توسعه یافته است تا توانمندسازی والدین، سرپرستان و کودکان ایجاد گرددهدف از این مطالعه، دریافت نظر و تجربهٔ شما در مورد پیام رسان بحران `"CrisisText`" است تا اطمینان حاصل گردد که این چت‌بات می‌تواند به خانواده‌هایی مانند شما کمک نماید.", 2) | Out-Null
$d.Content.Find.Execute("Before you decide if you’d like to join, it’s important for you to know why we’re doing this research and what it involves. You can read through this Participant Information Sheet. ", $true, $false, $false, $false, $false, $true, 1, $false, "Before you decide if you’d like to join, it’s important for you to know why we’re doing this research and what it involves. شما می‌توانید این ورقۀ معلومات اشتراک‌کننده را مطالعه نمایید. ", 2) | Out-Null
$d.Content.Find.Execute("If you have any questions about the chatbot or if something isn’t clear, please email the study team at ", $true, $false, $false, $false, $false, $true, 1, $false, "اگر در مورد چت‌بات پرسشی دارید یا موضوع برایتان روشن نیست، لطفاً با تیم تحقیق از طریق ایمیل ", 2) | Out-Null
$d.Content.Find.Execute(" or message us on WhatsApp at +27 79 762 3598. We’re here to help you! ", $true, $false, $false, $false, $false, $true, 1, $false, " یا پیام واتس‌اپ به شمارهٔ +27 79 762 3598 تماس بگیرید. ما برای راهنمایی و کمک در کنار شما هستیم! ", 2) | Out-Null
$d.Content.Find.Execute("Who can join?", $true, $false, $false, $false, $false, $true, 1, $false, "چی کسانی میتوانند اشتراک کنند؟", 2) | Out-Null
$d.Content.Find.Execute("To be part of the study, you need to be 18 years or older, be the parent or caregiver of a child under 18 years old, and live in a country with a participating World Vision office. You also need to agree to take part in the study. ", $true, $false, $false, $false, $false, $true, 1, $false, "برای شامل‌شدن در این ارزیابی، لازم است که ۱۸ ساله یا بزرگ‌تر باشید، والد یا مراقب یک کودک زیر ۱۸ سال بوده و در کشوری زندگی کنید که دفتر ورلد ویژن در آن اشتراک دارد. و شما برای اشتراک درین ارزیابی نیاز است تا موافقت کنید. ", 2) | Out-Null
$d.Content.Find.Execute("Do I have to join?", $true, $false, $false, $false, $false, $true, 1, $false, "آیا لازم است من اشتراک کنم؟", 2) | Out-Null
$d.Content.Find.Execute("No, it's up to you if you want to join or not. If you don't want to participate, nothing bad will happen to you or your family. If you do choose to join, you can stop at any time. If you want to stop getting messages, you can type `"STOP MESSAGES`". Additionally, if you participate but don’t want to answer some questions that the chatbot asks, you can simply skip any questions. You can still get the chatbot messages even if you don't answer the questions.", $true, $false, $false, $false, $false, $true, 1, $false, "نخیر، این بستگی به شما دارد که اشتراک می کنید یا خیر. اگر نمی‌خواهید اشتراک کنید، هیچ‌گونه پیامد بدی برای شما یا خانواده‌تان به‌وجود نخواهد آمد. اگر تصمیم به اشتراک گرفتید، می‌توانید در هر زمان که خواستید، اشتراک خود را متوقف سازید. اگر می‌خواهید دریافت پیام‌ها را متوقف کنید، می‌توانید عبارت «توقف پیام‌ها» را بنویسید. همچنین، اگر در این ارزیابی اشتراک کنید اما نخواهید به برخی پرسش‌هایی که چت‌بات می‌پرسد پاسخ دهید، می‌توانید آن پرسش‌ها را به‌سادگی نادیده بگیرید. شما حتی اگر به پرسش‌ها پاسخ ندهید، باز هم می‌توانید پیام‌های چت‌بات را دریافت کنید.", 2) | Out-Null
$d.Content.Find.Execute("What happens if I join?", $true, $false, $false, $false, $false, $true, 1, $false, "اگر من اشتراک کنم چی می شود؟", 2) | Out-Null
$d.Content.Find.Execute("If you decide to join, you’ll need to read the consent form below and answer “Yes” on WhatsApp to the question, “I have read and understand the information, and I give consent to participate in the study”. ", $true, $false, $false, $false, $false, $true, 1, $false, "اگر تصمیم به اشتراک گرفتید، نیاز است فرم رضایت‌نامهٔ زیر را مطالعه کنید و در واتس‌اپ به این پرسش «من معلومات را خوانده و درک کرده‌ام و برای اشتراک در ارزیابی رضایت می‌دهم» با نوشتن «بلی» پاسخ دهید. ", 2) | Out-Null
$d.Content.Find.Execute("The study team will then send you a survey through WhatsApp with about 8 questions. ", $true, $false, $false, $false, $false, $true, 1, $false, "سپس تیم ارزیابی از طریق واتس‌اپ یک پرسش‌نامه با حدود ۸ پرسش برایتان ارسال خواهد کرد. ", 2) | Out-Null
$d.Content.Find.Execute("The survey asks questions related to parenting and your well-being.", $true, $false, $false, $false, $false, $true, 1, $false, "این پرسش‌ها مرتبط به والدگری و سلامت روان و بهبود وضعیت شما می‌باشد.", 2) | Out-Null
$d.Content.Find.Execute("After this, we’ll begin the CrisisText programme. This programme includes sessions with tips aiming ", $true, $false, $false, $false, $false, $true, 1, $false, "پس از آن، برنامهٔ «پیام رسان بحران» را آغاز خواهیم کرد. این برنامه شامل جلساتی با توصیه های است که", 2) | Out-Null
$d.Content.Find.Execute("provide encouragement and actionable tips to:", $true, $false, $false, $false, $false, $true, 1, $false, " هدف آن تشویق و ارائهٔ رهنمودهای عملی برای موارد ذیل می‌باشد:", 2) | Out-Null
$d.Content.Find.Execute("1) Help parents heal from depression, anxiety, and trauma; ", $true, $false, $false, $false, $false, $true, 1, $false, "1) کمک به والدین برای بهبود افسردگی، اضطراب و صدمات روانی؛ ", 2) | Out-Null
$d.Content.Find.Execute("2) Improve parenting practices to keep children safe and healthy amidst crisis ", $true, $false, $false, $false, $false, $true, 1, $false, "2) ارتقای روش‌های تربیهٔ طفل برای حفظ صحت و مصونیت اطفال در شرایط بحرانی می‌باشد ", 2) | Out-Null
$d.Content.Find.Execute("We’ll send you a survey with questions again after you participate in CrisisText, and again one month later.", $true, $false, $false, $false, $false, $true, 1, $false, "یک پرسشنامه پس از پایان برنامه و پرسشنامهٔ دیگر یک ماه بعد برای شما ارسال خواهد شد.", 2) | Out-Null
$d.Content.Find.Execute("Do I get anything for joining? ", $true, $false, $false, $false, $false, $true, 1, $false, "آیا برای اشتراک چیزی دریافت میکنم؟ ", 2) | Out-Null
$d.Content.Find.Execute("We hope the CrisisText programme will provide helpful tips to build strength, hope, and encouragement for you and your children!", $true, $false, $false, $false, $false, $true, 1, $false, "امیدواریم برنامهٔ «پیام رسان بحران» نکات سودمندی را برای تقویت توانمندی، امید و دلگرمی برای شما و اطفال تان فراهم سازد!", 2) | Out-Null
$d.Content.Find.Execute("What happens to my information if I join?", $true, $false, $false, $false, $false, $true, 1, $false, "اگر اشتراک کنم، با معلوماتم چه می‌شود؟", 2) | Out-Null
$d.Content.Find.Execute("We only collect what’s needed for the study and store it securely. We will never ask for your name, and we do not ask for other identifying information such as date of birth. We also do not store your phone number. We only collect basic demographics (age, sex, number of children, and country).", $true, $false, $false, $false, $false, $true, 1, $false, "ما تنها معلومات لازم برای این ارزیابی را جمع‌آوری می‌کنیم و آن را به‌گونهٔ مصئون و امن نگهداری می‌نماییم. ما هرگز نام شما را نمی‌پرسیم و همچنان معلومات شناسایی‌کنندهٔ دیگر مانند تاریخ تولد را نیز درخواست نمی‌کنیم. ما همچنین شماره تماس شما را ذخیره نمی کنیم. ما تنها معلومات ابتدایی را جمع‌آوری می‌کنیم (سن، جنسیت، تعداد کودکان و کشور).", 2) | Out-Null
$d.Content.Find.Execute("Your information, including the answers you give during the surveys and data on how much of the program you complete, will be kept safe on secure servers connected with IDEMS, PLH, WV, and University of Oxford. All data will be kept for five years after the study. Ethics committees and monitors may check the information. We cannot remove or change your information, even if you stop participating; the reason for this is that we do not collect your identifying information (such as your name), and therefore, we won’t know which information is yours to remove or make changes. After the study, we may share the information with other researchers, but it will never be possible to know who took part.  ", $true, $false, $false, $false, $false, $true, 1, $false, "معلومات شما، شامل پاسخ‌هایی که در جریان پرسش‌نامه‌ها می‌دهید و همچنان معلومات مربوط به میزان تکمیل برنامه از سوی شما، در سرورهای امن مربوط به مؤسسات IDEMS، PLH، ورلد ویژن جهانی WV و دانشگاه آکسفورد محفوظ و نگهداری خواهد شد. تمام معلومات برای مدت پنج سال پس از پایان ارزیابی نگهداری خواهد شد. ناضران ممکن است معلومات شما را بررسی کنند. ما نمی‌توانیم معلومات شما را حذف یا تغییر دهیم، حتی اگر اشتراک‌تان را متوقف کنید؛ دلیل این موضوع آن است که ما معلومات هویتی شما (مانند نام‌تان) را جمع‌آوری نمی‌کنیم، بنابراین نمی‌دانیم کدام معلومات مربوط به شما است تا آن را حذف کرده یا تغییر دهیم. پس از پایان ارزیابی، ممکن است معلومات را با پژوهشگران دیگر شریک سازیم، اما هرگز امکان شناسایی افراد اشتراک‌کننده وجود نخواهد داشت.  ", 2) | Out-Null
$d.Content.Find.Execute("What happens to the research results?", $true, $false, $false, $false, $false, $true, 1, $false, "نتایج این ارزیابی چی می شود؟", 2) | Out-Null
$d.Content.Find.Execute("Your participation and what you tell us will help us understand how to support families like yours. We plan to share the results in journals, policy briefs and conferences so others can learn from this study too. When results are shared, it will not be possible to know who took part. We will never ask for your name.", $true, $false, $false, $false, $false, $true, 1, $false, "اشتراک شما و معلوماتی که در اختیار ما قرار می‌دهید، به ما کمک می‌کند تا دریابیم چگونه می‌توانیم از خانواده‌هایی مانند خانوادهٔ شما حمایت کنیم. ما قصد داریم نتایج این ارزیابی را در ژورنال‌ها، یادداشت‌های پالیسی و کنفرانس‌ها به‌اشتراک بگذاریم تا دیگران نیز بتوانند از یافته‌های این مطالعه بهره‌مند شوند. وقتی نتایج منتشر شود، شناسایی افراد اشتراک‌کننده ممکن نخواهد بود. ما هرگز نام شما را پرسان نخواهیم کرد.", 2) | Out-Null
$d.Content.Find.Execute("Who is in the study team?", $true, $false, $false, $false, $false, $true, 1, $false, "کی ها در تیم ارزیابی هستند؟", 2) | Out-Null
$d.Content.Find.Execute("The researchers of this study are Dr Jamie Lachman (Universities of Cape Town and Oxford) and Sydney Tucker (University of Oxford).", $true, $false, $false, $false, $false, $true, 1, $false, "پژوهشگران این ارزیابی داکتر جیمی لاچمن (دانشگاه‌های کیپ‌تاون و آکسفورد) و سیدنی تاکر (دانشگاه آکسفورد) می‌باشند.", 2) | Out-Null
$d.Content.Find.Execute("Are there any risks in joining?   ", $true, $false, $false, $false, $false, $true, 1, $false, "آیا اشتراک در این ارزیابی با کدام خطر همراه است؟   ", 2) | Out-Null
$d.Content.Find.Execute("We don’t expect any risks to you if you join this study. We hope the CrisisText programme will provide helpful tips to build strength, hope, and encouragement for you and your children!", $true, $false, $false, $false, $false, $true, 1, $false, "ما انتظار نداریم که اشتراک شما در این ارزیابی با هیچ‌گونه خطری همراه باشد. امیدواریم برنامهٔ «CrisisText» نکات سودمندی برای تقویت توانمندی، امید و دلگرمی برای شما و اطفال تان فراهم سازد!", 2) | Out-Null
$d.Content.Find.Execute("If any questions make you uncomfortable, you don’t have to answer them. If you become upset when using the programme, you can type “HELP” to receive troubleshooting messages, including resources in your local context.", $true, $false, $false, $false, $false, $true, 1, $false, "اگر پاسخ‌دادن به پرسشی برایتان ناراحت‌کننده باشد، لازم نیست به آن پاسخ دهید. اگر هنگام استفاده از برنامه احساس ناراحتی کردید، می‌توانید عبارت «کمک» را بنویسید تا پیام‌های راهنمایی برای رفع مشکل، شامل منابع مربوط به محل زندگی‌تان، دریافت نمایید.", 2) | Out-Null
$d.Content.Find.Execute("Additionally, remember, you can stop participating anytime, without giving a reason. ", $true, $false, $false, $false, $false, $true, 1, $false, "همچنین، به‌خاطر داشته باشید که می‌توانید در هر زمان بدون ارائهٔ دلیل، اشتراک خود را متوقف سازید. ", 2) | Out-Null
$d.Content.Find.Execute("Who pays for the study?", $true, $false, $false, $false, $false, $true, 1, $false, "هزینهٔ این ارزیابی را چه‌کسی پرداخت می‌کند؟", 2) | Out-Null
$d.Content.Find.Execute("This study is funded through Parenting for Lifelong Health, World Vision, and the Global Parenting Initiative, funded by the LEGO Foundation (CVR00940), Oak Foundation, the World Childhood Foundation (16191), The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund (ES/S008101/1). ", $true, $false, $false, $false, $false, $true, 1, $false, "این پژوهش توسط برنامهٔ «والدگری برای زندگی سالم» (Parenting for Lifelong Health)، ورلد ویژن (World Vision) و ابتکار جهانی والدگری Global Parenting Initiative تمویل می‌شود. منابع مالی آن از بنیاد لیگو (LEGO Foundation) با شمارهٔ CVR00940، بنیاد اوک Oak Foundatio، بنیاد جهانی کودک World Childhood Foundation با شمارهٔ 16191، برنامهٔ Human Safety Net، و صندوق پژوهشی چالش‌های جهانی وزارت پژوهش و نوآوری بریتانیا UK Research and Innovation Global Challenges Research Fund با شمارهٔ ES/S008101/1 فراهم گردیده است. ", 2) | Out-Null

Write-Host "Replacements applied"
